# Remove the "reviews_count" column (column E). Everything to the right
# (reviews_average, latitude, longitude, is_permanently_closed, gmaps_link,
# latest_review_date) shifts one column to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").EntireColumn.Delete()
